# Slide 12 ("Water Utility: Relationships & Cardinalities") contains two
# shapes: the title textbox (Shape 117) and the ER-diagram picture
# ("Picture 3"). The edit repositions the picture:
#   a:off x="713232" y="662249"  ->  a:off x="240386" y="677489"   (EMU)
# Size (a:ext) is unchanged.
#
# PowerPoint's COM object model reports/accepts Left/Top in points
# (1 pt = 12700 EMU), so convert the target EMU values to points. A
# sub-EMU nudge (+0.5 EMU worth of points) is added before dividing to
# offset float32 rounding in the points -> EMU round trip, so the saved
# XML lands exactly on the target integer EMU values.

$EmuPerPoint   = 12700.0
$TargetLeftEmu = 240386
$TargetTopEmu  = 677489

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(2)   # "Picture 3"

$sh.Left = ($TargetLeftEmu + 0.5) / $EmuPerPoint
$sh.Top  = ($TargetTopEmu  + 0.5) / $EmuPerPoint
